$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new top data row for "2022-Q3"
#    above the existing "2022-Q1" row, pushing all following rows down by
#    one and extending the last row (2020-Q4) accordingly.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Push existing data rows (old rows 2-7) down to make room for the new row.
$summary.Rows.Item(2).Insert()

# Column A on the summary sheet carries the bold/bordered "index" style;
# grab it from the row directly below (the old "2022-Q1" row, now row 3)
# and drop the bold/border formatting that Insert() copied into B:D from
# the header row above.
$summary.Cells.Item(3,1).Copy()
$summary.Cells.Item(2,1).PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 4
$summary.Cells.Item(2,4).Value = 0.08

# ---------------------------------------------------------------------------
# 2) Add the new "2022-Q3" detail sheet, positioned right before the
#    "2022-Q1" sheet. Duplicating "2022-Q1" keeps the page setup / outline
#    properties / header styling identical to its siblings.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1Index = $q1.Index
$q1.Copy($q1)
$q3 = $wb.Worksheets.Item($q1Index)
$q3.Name = "2022-Q3"

# Replicate the formatted (but empty) row 2 down through row 5 so every
# data row gets the same "index column is bold/bordered, rest are plain"
# styling as row 2.
$q3.Range("A2:H2").Copy()
$q3.Range("A3:H5").PasteSpecial(-4122)

# Force columns B:G to text so numeric-looking values (fund codes,
# percentages, etc.) are stored as strings, matching the source data.
$q3.Range("B2:G5").NumberFormat = "@"

$rows = @(
    @(0, "159851", "华宝中证金融科技主题ETF", "1.94", "98.27", "2.50", "0.0485", 10),
    @(1, "159628", "万家国证2000ETF", "2.90", "97.72", "0.45", "0.0130", 9),
    @(2, "516100", "华夏中证金融科技主题ETF", "0.51", "96.79", "2.47", "0.0126", 10),
    @(3, "516860", "博时中证金融科技主题ETF", "0.34", "98.57", "2.51", "0.0085", 10)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r,1).Value = $row[0]
    $q3.Cells.Item($r,2).Value = $row[1]
    $q3.Cells.Item($r,3).Value = $row[2]
    $q3.Cells.Item($r,4).Value = $row[3]
    $q3.Cells.Item($r,5).Value = $row[4]
    $q3.Cells.Item($r,6).Value = $row[5]
    $q3.Cells.Item($r,7).Value = $row[6]
    $q3.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

# Drop the temporary text NumberFormat now that the values are locked in as
# strings, so the cells end up with no explicit style (matching siblings).
$q3.Range("B2:G5").ClearFormats()
